# Replace the single-letter labels on the "Oval" shapes across every slide:
#   A -> X
#   B -> Z
#   C -> Y
# (the three ovals used throughout the deck as step markers)

$p = $ppt.ActivePresentation

$map = @{ "A" = "X"; "B" = "Z"; "C" = "Y" }

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shape = $s.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            $old = $tr.Text
            if ($map.ContainsKey($old)) {
                $tr.Text = $map[$old]
            }
        }
    }
}
